# Update the Protocol_Name value for the test data row from
# "IEC 61850 Ed1" to "IEC 61850 Ed2", and move the active
# selection to the edited cell (G2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "IEC 61850 Ed2"

$ws.Range("G2").Select()
